# Auto-generated script applying 2022-11-23 violent crime data update
# Updates 2022 (column I) totals across the Citywide Totals, By Neighborhood,
# and individual neighborhood worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("I2").Value = 6567  # was 6548
$ws.Range("I3").Value = 6854  # was 6832
$ws.Range("I4").Value = 1575  # was 1572
$ws.Range("I6").Value = 7883  # was 7853
$ws.Range("I7").Value = 23516  # was 23442

$ws = $wb.Worksheets.Item(2)
$ws.Range("I2").Value = 183  # was 182
$ws.Range("I4").Value = 96  # was 95
$ws.Range("I5").Value = 70  # was 69
$ws.Range("I7").Value = 745  # was 743
$ws.Range("I8").Value = 1406  # was 1399
$ws.Range("I9").Value = 120  # was 119
$ws.Range("I10").Value = 169  # was 168
$ws.Range("I12").Value = 57  # was 56
$ws.Range("I15").Value = 274  # was 272
$ws.Range("I18").Value = 182  # was 180
$ws.Range("I19").Value = 665  # was 664
$ws.Range("I20").Value = 584  # was 581
$ws.Range("I23").Value = 230  # was 229
$ws.Range("I29").Value = 1421  # was 1414
$ws.Range("I33").Value = 1055  # was 1049
$ws.Range("I36").Value = 322  # was 321
$ws.Range("I37").Value = 743  # was 738
$ws.Range("I42").Value = 854  # was 851
$ws.Range("I45").Value = 47  # was 46
$ws.Range("I48").Value = 303  # was 302
$ws.Range("I52").Value = 516  # was 514
$ws.Range("I54").Value = 475  # was 474
$ws.Range("I55").Value = 268  # was 267
$ws.Range("I57").Value = 93  # was 91
$ws.Range("I60").Value = 132  # was 131
$ws.Range("I63").Value = 72  # was 78
$ws.Range("I67").Value = 901  # was 898
$ws.Range("I70").Value = 37  # was 36
$ws.Range("I73").Value = 215  # was 214
$ws.Range("I74").Value = 39  # was 38
$ws.Range("I77").Value = 142  # was 140
$ws.Range("I78").Value = 318  # was 317
$ws.Range("I79").Value = 668  # was 666
$ws.Range("I83").Value = 510  # was 509
$ws.Range("I85").Value = 1056  # was 1054
$ws.Range("I87").Value = 54  # was 53
$ws.Range("I88").Value = 220  # was 218
$ws.Range("I89").Value = 281  # was 279
$ws.Range("I90").Value = 304  # was 303
$ws.Range("I91").Value = 250  # was 249
$ws.Range("I94").Value = 238  # was 237
$ws.Range("I96").Value = 263  # was 261
$ws.Range("I99").Value = 416  # was 415
$ws.Range("I101").Value = 23516  # was 23442

$ws = $wb.Worksheets.Item(3)
$ws.Range("I2").Value = 299  # was 298
$ws.Range("I6").Value = 272  # was 271
$ws.Range("I7").Value = 1056  # was 1054

$ws = $wb.Worksheets.Item(5)
$ws.Range("I2").Value = 137  # was 136
$ws.Range("I6").Value = 148  # was 147
$ws.Range("I7").Value = 516  # was 514

$ws = $wb.Worksheets.Item(7)
$ws.Range("I2").Value = 417  # was 416
$ws.Range("I3").Value = 405  # was 404
$ws.Range("I4").Value = 90  # was 88
$ws.Range("I6").Value = 452  # was 449
$ws.Range("I7").Value = 1406  # was 1399

$ws = $wb.Worksheets.Item(9)
$ws.Range("I2").Value = 244  # was 243
$ws.Range("I6").Value = 200  # was 199
$ws.Range("I7").Value = 745  # was 743

$ws = $wb.Worksheets.Item(10)
$ws.Range("I2").Value = 68  # was 67
$ws.Range("I3").Value = 69  # was 68
$ws.Range("I7").Value = 281  # was 279

$ws = $wb.Worksheets.Item(11)
$ws.Range("I6").Value = 99  # was 97
$ws.Range("I7").Value = 263  # was 261

$ws = $wb.Worksheets.Item(14)
$ws.Range("I2").Value = 221  # was 219
$ws.Range("I3").Value = 243  # was 240
$ws.Range("I7").Value = 743  # was 738

$ws = $wb.Worksheets.Item(15)
$ws.Range("I6").Value = 107  # was 106
$ws.Range("I7").Value = 416  # was 415

$ws = $wb.Worksheets.Item(16)
$ws.Range("I3").Value = 331  # was 330
$ws.Range("I6").Value = 274  # was 272
$ws.Range("I7").Value = 901  # was 898

$ws = $wb.Worksheets.Item(20)
$ws.Range("I3").Value = 185  # was 184
$ws.Range("I7").Value = 510  # was 509

$ws = $wb.Worksheets.Item(22)
$ws.Range("I2").Value = 238  # was 237
$ws.Range("I3").Value = 388  # was 384
$ws.Range("I6").Value = 338  # was 337
$ws.Range("I7").Value = 1055  # was 1049

$ws = $wb.Worksheets.Item(24)
$ws.Range("I6").Value = 230  # was 229
$ws.Range("I7").Value = 475  # was 474

$ws = $wb.Worksheets.Item(25)
$ws.Range("I3").Value = 491  # was 488
$ws.Range("I4").Value = 72  # was 70
$ws.Range("I6").Value = 394  # was 392
$ws.Range("I7").Value = 1421  # was 1414

$ws = $wb.Worksheets.Item(26)
$ws.Range("I2").Value = 219  # was 218
$ws.Range("I7").Value = 665  # was 664

$ws = $wb.Worksheets.Item(28)
$ws.Range("I6").Value = 155  # was 154
$ws.Range("I7").Value = 303  # was 302

$ws = $wb.Worksheets.Item(32)
$ws.Range("I6").Value = 307  # was 304
$ws.Range("I7").Value = 854  # was 851

$ws = $wb.Worksheets.Item(34)
$ws.Range("I2").Value = 52  # was 51
$ws.Range("I7").Value = 169  # was 168

$ws = $wb.Worksheets.Item(35)
$ws.Range("I6").Value = 116  # was 115
$ws.Range("I7").Value = 318  # was 317

$ws = $wb.Worksheets.Item(36)
$ws.Range("I2").Value = 81  # was 80
$ws.Range("I7").Value = 268  # was 267

$ws = $wb.Worksheets.Item(39)
$ws.Range("I3").Value = 81  # was 80
$ws.Range("I7").Value = 230  # was 229

$ws = $wb.Worksheets.Item(40)
$ws.Range("I6").Value = 68  # was 67
$ws.Range("I7").Value = 250  # was 249

$ws = $wb.Worksheets.Item(42)
$ws.Range("I2").Value = 192  # was 190
$ws.Range("I7").Value = 668  # was 666

$ws = $wb.Worksheets.Item(44)
$ws.Range("I2").Value = 163  # was 161
$ws.Range("I3").Value = 165  # was 164
$ws.Range("I7").Value = 584  # was 581

$ws = $wb.Worksheets.Item(45)
$ws.Range("I3").Value = 41  # was 40
$ws.Range("I6").Value = 85  # was 84
$ws.Range("I7").Value = 182  # was 180

$ws = $wb.Worksheets.Item(47)
$ws.Range("I6").Value = 101  # was 100
$ws.Range("I7").Value = 322  # was 321

$ws = $wb.Worksheets.Item(51)
$ws.Range("I2").Value = 43  # was 44
$ws.Range("I6").Value = 138  # was 136
$ws.Range("I7").Value = 238  # was 237

$ws = $wb.Worksheets.Item(54)
$ws.Range("I4").Value = 16  # was 15
$ws.Range("I6").Value = 104  # was 103
$ws.Range("I7").Value = 274  # was 272

$ws = $wb.Worksheets.Item(61)
$ws.Range("I3").Value = 41  # was 40
$ws.Range("I7").Value = 120  # was 119

$ws = $wb.Worksheets.Item(62)
$ws.Range("I4").Value = 21  # was 20
$ws.Range("I7").Value = 215  # was 214

$ws = $wb.Worksheets.Item(64)
$ws.Range("I2").Value = 67  # was 66
$ws.Range("I7").Value = 183  # was 182

$ws = $wb.Worksheets.Item(67)
$ws.Range("I4").Value = 8  # was 7
$ws.Range("I7").Value = 37  # was 36

$ws = $wb.Worksheets.Item(68)
$ws.Range("I3").Value = 75  # was 74
$ws.Range("I6").Value = 69  # was 68
$ws.Range("I7").Value = 220  # was 218

$ws = $wb.Worksheets.Item(70)
$ws.Range("I3").Value = 19  # was 18
$ws.Range("I7").Value = 70  # was 69

$ws = $wb.Worksheets.Item(74)
$ws.Range("I6").Value = 107  # was 106
$ws.Range("I7").Value = 304  # was 303

$ws = $wb.Worksheets.Item(77)
$ws.Range("I2").Value = 34  # was 33
$ws.Range("I6").Value = 24  # was 23
$ws.Range("I7").Value = 93  # was 91

$ws = $wb.Worksheets.Item(78)
$ws.Range("I6").Value = 41  # was 40
$ws.Range("I7").Value = 132  # was 131

$ws = $wb.Worksheets.Item(84)
$ws.Range("I2").Value = 46  # was 45
$ws.Range("I6").Value = 34  # was 33
$ws.Range("I7").Value = 142  # was 140

$ws = $wb.Worksheets.Item(85)
$ws.Range("I2").Value = 18  # was 17
$ws.Range("I7").Value = 47  # was 46

$ws = $wb.Worksheets.Item(90)
$ws.Range("I2").Value = 35  # was 34
$ws.Range("I7").Value = 96  # was 95

$ws = $wb.Worksheets.Item(91)
$ws.Range("I6").Value = 30  # was 29
$ws.Range("I7").Value = 57  # was 56

$ws = $wb.Worksheets.Item(92)
$ws.Range("I3").Value = 10  # was 9
$ws.Range("I7").Value = 54  # was 53

$ws = $wb.Worksheets.Item(95)
$ws.Range("I2").Value = 10  # was 9
$ws.Range("I7").Value = 39  # was 38
